$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same set applied to every data row, columns B:Q)
$values = @(
    0.9999969141016266,
    0.9990169585194482,
    0.9999995420867596,
    0.9999988375479447,
    0.9999991128027239,
    0.00000288055067624043,
    0.0009176260715613002,
    0.0000003346038930132623,
    0.000001344763944861583,
    0.0000008396839189374227,
    0.00009999963621106398,
    0.00169721851163615,
    0.9999753128130129,
    0.001769472617403162,
    67.51505815013246,
    93.11145047236465
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
